# Add a "*" suffix filter marker to the gene name cells in columns J and V
# of the Combined_Genes sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined_Genes")

for ($row = 1; $row -le 11; $row++) {
    $cell = $ws.Range("J$row")
    $cell.Value = $cell.Value() + "*"
}

for ($row = 1; $row -le 13; $row++) {
    if ($row -eq 2) {
        continue
    }
    $cell = $ws.Range("V$row")
    $cell.Value = $cell.Value() + "*"
}
